# Auto update stock data
# Updates the "as of" date (column A) and EBITDA figure (column B) for each
# company's latest row. Values in this sheet are stored as plain text
# (inline/shared strings), so we force text interpretation via NumberFormat
# "@" before assigning, then restore the cell's original Style so no new
# formatting/style footprint is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($worksheet, $row, $col, $text) {
    $cell = $worksheet.Cells.Item($row, $col)
    $originalStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $originalStyle
}

# row -> (new date, new EBITDA value or $null if unchanged)
$updates = @{
    2  = @("2025/11/13", "4.89")
    8  = @("2025/11/13", "7.61")
    14 = @("2025/11/13", "2.89")
    20 = @("2025/11/13", "12.24")
    26 = @("2025/11/13", "9.92")
    32 = @("2025/11/13", "24.84")
    38 = @("2025/11/13", $null)
    44 = @("2025/11/13", "11.19")
    50 = @("2025/11/13", "11.69")
    56 = @("2025/11/13", "34.97")
    62 = @("2025/11/13", "11.62")
    68 = @("2025/11/13", "13.09")
    74 = @("2025/11/13", "15.82")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $newDate = $vals[0]
    $newEbitda = $vals[1]

    Set-TextCell $ws $row 1 $newDate

    if ($newEbitda -ne $null) {
        Set-TextCell $ws $row 2 $newEbitda
    }
}
